$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AD2").Value = 17
$ws.Range("AS2").Value = 15
$ws.Range("BF2").Value = "2013-01-17"
$ws.Range("AD3").Value = 17
$ws.Range("AS3").Value = 17
$ws.Range("BB3").Value = 20
$ws.Range("BF3").Value = "2013-01-17"
$ws.Range("AD4").Value = 9
$ws.Range("AH4").Value = 8
$ws.Range("AK4").Value = 18
$ws.Range("AS4").Value = 27
$ws.Range("BA4").Value = 6
$ws.Range("BF4").Value = "2013-01-17"
$ws.Range("AD5").Value = 17
$ws.Range("BB5").Value = 24
$ws.Range("BF5").Value = "2013-01-17"
$ws.Range("AD6").Value = 23
$ws.Range("AJ6").Value = 25
$ws.Range("AR6").Value = 11
$ws.Range("BF6").Value = "2013-01-17"
$ws.Range("AQ7").Value = 21
$ws.Range("BB7").Value = 22
$ws.Range("BF7").Value = "2013-01-17"
$ws.Range("AD8").Value = 4
$ws.Range("AN8").Value = 13
$ws.Range("AT8").Value = 17
$ws.Range("AW8").Value = 16
$ws.Range("BF8").Value = "2013-01-17"
$ws.Range("AE9").Value = 5
$ws.Range("AH9").Value = 22
$ws.Range("AP9").Value = 5
$ws.Range("BF9").Value = "2013-01-17"
$ws.Range("D10").Value = 38
$ws.Range("F10").Value = 24
$ws.Range("G10").Value = 0.368
$ws.Range("H10").Value = 48.7
$ws.Range("N10").Value = 0.374
$ws.Range("O10").Value = 16.8
$ws.Range("P10").Value = 23.6
$ws.Range("Q10").Value = 0.714
$ws.Range("S10").Value = 31.1
$ws.Range("T10").Value = 43.6
$ws.Range("AB10").Value = 95
$ws.Range("AC10").Value = -1.2
$ws.Range("AD10").Value = 17
$ws.Range("AF10").Value = 22
$ws.Range("AG10").Value = 23
$ws.Range("AH10").Value = 6
$ws.Range("AK10").Value = 17
$ws.Range("AO10").Value = 16
$ws.Range("AT10").Value = 8
$ws.Range("AZ10").Value = 17
$ws.Range("BA10").Value = 11
$ws.Range("BB10").Value = 21
$ws.Range("BC10").Value = 17
$ws.Range("BF10").Value = "2013-01-17"
$ws.Range("AD11").Value = 23
$ws.Range("AK11").Value = 9
$ws.Range("AN11").Value = 4
$ws.Range("AO11").Value = 17
$ws.Range("AY11").Value = 17
$ws.Range("BF11").Value = "2013-01-17"
$ws.Range("AD12").Value = 4
$ws.Range("AG12").Value = 14
$ws.Range("AN12").Value = 15
$ws.Range("AQ12").Value = 16
$ws.Range("AX12").Value = 28
$ws.Range("AY12").Value = 23
$ws.Range("AZ12").Value = 11
$ws.Range("BC12").Value = 12
$ws.Range("BF12").Value = "2013-01-17"
$ws.Range("AD13").Value = 4
$ws.Range("AE13").Value = 5
$ws.Range("AQ13").Value = 24
$ws.Range("BF13").Value = "2013-01-17"
$ws.Range("D14").Value = 39
$ws.Range("E14").Value = 30
$ws.Range("G14").Value = 0.769
$ws.Range("I14").Value = 38.7
$ws.Range("J14").Value = 80.90000000000001
$ws.Range("M14").Value = 20.6
$ws.Range("N14").Value = 0.355
$ws.Range("O14").Value = 17.5
$ws.Range("P14").Value = 24.4
$ws.Range("Q14").Value = 0.718
$ws.Range("S14").Value = 30.3
$ws.Range("T14").Value = 41.7
$ws.Range("U14").Value = 23.8
$ws.Range("V14").Value = 14.4
$ws.Range("X14").Value = 6.3
$ws.Range("AA14").Value = 21.4
$ws.Range("AB14").Value = 102.1
$ws.Range("AC14").Value = 8.9
$ws.Range("AD14").Value = 9
$ws.Range("AE14").Value = 2
$ws.Range("AJ14").Value = 24
$ws.Range("AR14").Value = 17
$ws.Range("AS14").Value = 18
$ws.Range("AT14").Value = 19
$ws.Range("AV14").Value = 10
$ws.Range("BA14").Value = 5
$ws.Range("BF14").Value = "2013-01-17"
$ws.Range("D15").Value = 38
$ws.Range("F15").Value = 21
$ws.Range("G15").Value = 0.447
$ws.Range("I15").Value = 37.4
$ws.Range("J15").Value = 81.90000000000001
$ws.Range("M15").Value = 25.3
$ws.Range("N15").Value = 0.355
$ws.Range("O15").Value = 19.7
$ws.Range("Q15").Value = 0.6919999999999999
$ws.Range("R15").Value = 12.5
$ws.Range("S15").Value = 32.8
$ws.Range("V15").Value = 15.3
$ws.Range("W15").Value = 7.6
$ws.Range("X15").Value = 5.7
$ws.Range("Y15").Value = 5
$ws.Range("Z15").Value = 19.2
$ws.Range("AA15").Value = 23.1
$ws.Range("AB15").Value = 103.4
$ws.Range("AC15").Value = 1.9
$ws.Range("AD15").Value = 17
$ws.Range("AJ15").Value = 18
$ws.Range("AK15").Value = 8
$ws.Range("AN15").Value = 16
$ws.Range("AR15").Value = 10
$ws.Range("AV15").Value = 26
$ws.Range("AW15").Value = 18
$ws.Range("AX15").Value = 9
$ws.Range("AZ15").Value = 7
$ws.Range("BC15").Value = 11
$ws.Range("BF15").Value = "2013-01-17"
$ws.Range("AD16").Value = 23
$ws.Range("AE16").Value = 5
$ws.Range("AG16").Value = 5
$ws.Range("AZ16").Value = 18
$ws.Range("BF16").Value = "2013-01-17"
$ws.Range("D17").Value = 37
$ws.Range("E17").Value = 25
$ws.Range("G17").Value = 0.676
$ws.Range("J17").Value = 78.2
$ws.Range("L17").Value = 8.6
$ws.Range("M17").Value = 22.1
$ws.Range("N17").Value = 0.391
$ws.Range("Q17").Value = 0.764
$ws.Range("R17").Value = 8.1
$ws.Range("T17").Value = 38.9
$ws.Range("V17").Value = 13.8
$ws.Range("W17").Value = 8.199999999999999
$ws.Range("X17").Value = 5.2
$ws.Range("Y17").Value = 3.2
$ws.Range("Z17").Value = 19.2
$ws.Range("AA17").Value = 19.9
$ws.Range("AC17").Value = 5.3
$ws.Range("AD17").Value = 23
$ws.Range("AH17").Value = 13
$ws.Range("AM17").Value = 7
$ws.Range("AN17").Value = 3
$ws.Range("AW17").Value = 12
$ws.Range("AX17").Value = 16
$ws.Range("AZ17").Value = 9
$ws.Range("BF17").Value = "2013-01-17"
$ws.Range("D18").Value = 37
$ws.Range("E18").Value = 19
$ws.Range("G18").Value = 0.514
$ws.Range("I18").Value = 37.1
$ws.Range("J18").Value = 85.90000000000001
$ws.Range("M18").Value = 18.1
$ws.Range("N18").Value = 0.336
$ws.Range("P18").Value = 21.4
$ws.Range("Q18").Value = 0.746
$ws.Range("S18").Value = 30.6
$ws.Range("T18").Value = 43.5
$ws.Range("U18").Value = 21.7
$ws.Range("X18").Value = 7.6
$ws.Range("AA18").Value = 19.5
$ws.Range("AC18").Value = -1.2
$ws.Range("AD18").Value = 23
$ws.Range("AE18").Value = 17
$ws.Range("AG18").Value = 16
$ws.Range("AH18").Value = 23
$ws.Range("AI18").Value = 13
$ws.Range("AO18").Value = 22
$ws.Range("AP18").Value = 19
$ws.Range("AQ18").Value = 19
$ws.Range("AS18").Value = 16
$ws.Range("AT18").Value = 9
$ws.Range("AU18").Value = 19
$ws.Range("AZ18").Value = 9
$ws.Range("BA18").Value = 20
$ws.Range("BC18").Value = 18
$ws.Range("BF18").Value = "2013-01-17"
$ws.Range("D19").Value = 35
$ws.Range("F19").Value = 19
$ws.Range("G19").Value = 0.457
$ws.Range("I19").Value = 35.3
$ws.Range("K19").Value = 0.429
$ws.Range("N19").Value = 0.297
$ws.Range("O19").Value = 19
$ws.Range("P19").Value = 25.8
$ws.Range("Q19").Value = 0.738
$ws.Range("S19").Value = 31.5
$ws.Range("T19").Value = 45.2
$ws.Range("U19").Value = 21.8
$ws.Range("V19").Value = 15.3
$ws.Range("AA19").Value = 23
$ws.Range("AB19").Value = 95.3
$ws.Range("AC19").Value = -1.6
$ws.Range("AD19").Value = 30
$ws.Range("AF19").Value = 15
$ws.Range("AO19").Value = 5
$ws.Range("AP19").Value = 4
$ws.Range("AQ19").Value = 23
$ws.Range("AU19").Value = 16
$ws.Range("AW19").Value = 15
$ws.Range("AY19").Value = 25
$ws.Range("AZ19").Value = 2
$ws.Range("BB19").Value = 19
$ws.Range("BF19").Value = "2013-01-17"
$ws.Range("AD20").Value = 9
$ws.Range("AK20").Value = 14
$ws.Range("AO20").Value = 26
$ws.Range("AS20").Value = 18
$ws.Range("AT20").Value = 18
$ws.Range("AV20").Value = 11
$ws.Range("BF20").Value = "2013-01-17"
$ws.Range("D21").Value = 37
$ws.Range("E21").Value = 24
$ws.Range("G21").Value = 0.649
$ws.Range("I21").Value = 37.2
$ws.Range("J21").Value = 83.7
$ws.Range("M21").Value = 28.8
$ws.Range("N21").Value = 0.388
$ws.Range("O21").Value = 15.9
$ws.Range("P21").Value = 21.2
$ws.Range("Q21").Value = 0.749
$ws.Range("S21").Value = 29.5
$ws.Range("T21").Value = 40.5
$ws.Range("U21").Value = 20.1
$ws.Range("V21").Value = 11.1
$ws.Range("AC21").Value = 4.7
$ws.Range("AD21").Value = 23
$ws.Range("AI21").Value = 12
$ws.Range("AJ21").Value = 6
$ws.Range("AK21").Value = 15
$ws.Range("AO21").Value = 23
$ws.Range("AP21").Value = 21
$ws.Range("AQ21").Value = 17
$ws.Range("AT21").Value = 27
$ws.Range("AU21").Value = 26
$ws.Range("AW21").Value = 8
$ws.Range("AX21").Value = 29
$ws.Range("BF21").Value = "2013-01-17"
$ws.Range("AD22").Value = 9
$ws.Range("AU22").Value = 18
$ws.Range("AW22").Value = 11
$ws.Range("BA22").Value = 10
$ws.Range("BF22").Value = "2013-01-17"
$ws.Range("AD23").Value = 17
$ws.Range("AH23").Value = 15
$ws.Range("AK23").Value = 10
$ws.Range("BB23").Value = 23
$ws.Range("BF23").Value = "2013-01-17"
$ws.Range("AD24").Value = 9
$ws.Range("AJ24").Value = 7
$ws.Range("AM24").Value = 23
$ws.Range("AS24").Value = 18
$ws.Range("AU24").Value = 17
$ws.Range("AW24").Value = 21
$ws.Range("BF24").Value = "2013-01-17"
$ws.Range("D25").Value = 40
$ws.Range("F25").Value = 27
$ws.Range("G25").Value = 0.325
$ws.Range("I25").Value = 37.5
$ws.Range("J25").Value = 84.3
$ws.Range("K25").Value = 0.445
$ws.Range("L25").Value = 6.1
$ws.Range("M25").Value = 18.3
$ws.Range("N25").Value = 0.333
$ws.Range("O25").Value = 14.3
$ws.Range("P25").Value = 19.3
$ws.Range("Q25").Value = 0.742
$ws.Range("R25").Value = 11.5
$ws.Range("S25").Value = 29.1
$ws.Range("T25").Value = 40.6
$ws.Range("U25").Value = 21.9
$ws.Range("V25").Value = 14
$ws.Range("W25").Value = 7.6
$ws.Range("Y25").Value = 5.1
$ws.Range("Z25").Value = 20.5
$ws.Range("AD25").Value = 4
$ws.Range("AH25").Value = 18
$ws.Range("AK25").Value = 16
$ws.Range("AM25").Value = 22
$ws.Range("AO25").Value = 27
$ws.Range("AS25").Value = 28
$ws.Range("AT25").Value = 25
$ws.Range("AW25").Value = 17
$ws.Range("AY25").Value = 16
$ws.Range("BF25").Value = "2013-01-17"
$ws.Range("AD26").Value = 9
$ws.Range("AP26").Value = 20
$ws.Range("AW26").Value = 19
$ws.Range("BF26").Value = "2013-01-17"
$ws.Range("AD27").Value = 9
$ws.Range("AT27").Value = 24
$ws.Range("AU27").Value = 27
$ws.Range("BF27").Value = "2013-01-17"
$ws.Range("AE28").Value = 2
$ws.Range("AV28").Value = 27
$ws.Range("AZ28").Value = 1
$ws.Range("BF28").Value = "2013-01-17"
$ws.Range("AD29").Value = 9
$ws.Range("AF29").Value = 25
$ws.Range("AG29").Value = 25
$ws.Range("AM29").Value = 8
$ws.Range("AS29").Value = 26
$ws.Range("AY29").Value = 18
$ws.Range("BF29").Value = "2013-01-17"
$ws.Range("AD30").Value = 4
$ws.Range("AG30").Value = 14
$ws.Range("AJ30").Value = 19
$ws.Range("AO30").Value = 6
$ws.Range("AW30").Value = 10
$ws.Range("AY30").Value = 23
$ws.Range("BF30").Value = "2013-01-17"
$ws.Range("AQ31").Value = 20
$ws.Range("AW31").Value = 20
$ws.Range("BF31").Value = "2013-01-17"
